$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the dates in A2:A32 forward by 31 days (2021-07 -> 2021-08 series)
for ($r = 2; $r -le 32; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $cur + 31
}

# Rows 33:63 held a duplicated copy of the July data (station "Балхаш I").
# Clear it out entirely: drop C/D (station name columns) completely and
# blank A/B/E back to bare, styled, empty cells.
$ws.Range("C33:D63").ClearContents()
$ws.Range("A33:A63").Value = ""
$ws.Range("B33:B63").Value = ""
$ws.Range("E33:E63").Value = ""

# Update the visible selection to A2:A32 (active cell A2)
$ws.Range("A2:A32").Select() | Out-Null

# Best-effort: nudge the window's vertical position (yWindow 6000 -> 7200 in the XML)
$win = $excel.ActiveWindow
$win.Top = 480
